# Generate Report for Handback
# This script mirrors a new handback pass: "a.md" and "b.md" swap roles
# (the row that used to describe a.md now describes b.md and vice versa),
# b.md's status flips to "not in sync" and a.md receives a fresh handback
# (new timestamps + handback name "TestHandback_201702211106").

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

$ov.Range("A2").Value = "b.md"
$ov.Range("B2").Value = "e2e\b.md"
$ov.Range("E2").Value = "Handed back: not in sync with en-US"
$ov.Range("F2").Value = "Handed back: not in sync with en-US"

$ov.Range("A3").Value = "a.md"
$ov.Range("B3").Value = "e2e\a.md"
$ov.Range("G3").Value = "2017-02-21 03:05:14"

# Hyperlinks on column B need their display text swapped too. Deleting via
# a range's Hyperlinks collection clears the sheet's hyperlinks so we
# recreate both, keeping them pointed at the same targets (and therefore
# the same relationship ids) as before.
$ov.Range("B2").Hyperlinks.Delete()
$ov.Hyperlinks.Add($ov.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test4/blob/5b8f33dc302b5a2aa99f42855abaa4d3b6b8492e/e2e/a.md", "", "", "e2e\b.md") | Out-Null
$ov.Hyperlinks.Add($ov.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test4/blob/5b8f33dc302b5a2aa99f42855abaa4d3b6b8492e/e2e/b.md", "", "", "e2e\a.md") | Out-Null

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("A2").Value = "b.md"
$zh.Range("C2").Value = "Handed back: not in sync with en-US"
$zh.Range("J2").Value = "b.md"

$zh.Range("A3").Value = "a.md"
$zh.Range("H3").Value = "2017-02-21 03:04:57"
$zh.Range("J3").Value = "a.md"
$zh.Range("L3").Value = "2017-02-21 03:06:28"
$zh.Range("M3").Value = "TestHandback_201702211106"

$zh.Range("A2").Hyperlinks.Delete()
$zh.Hyperlinks.Add($zh.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test4/blob/5b8f33dc302b5a2aa99f42855abaa4d3b6b8492e/e2e/a.md", "", "", "b.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("J2"), "https://github.com/OpenLocalizationTestOrg/ol-test4-zhcn/blob/23effdb786f34a4dbb604aff3a93b3ae491010c3/e2e/a.md", "", "", "b.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test4/blob/5b8f33dc302b5a2aa99f42855abaa4d3b6b8492e/e2e/b.md", "", "", "a.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("J3"), "https://github.com/OpenLocalizationTestOrg/ol-test4-zhcn/blob/23effdb786f34a4dbb604aff3a93b3ae491010c3/e2e/b.md", "", "", "a.md") | Out-Null

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("A2").Value = "b.md"
$de.Range("C2").Value = "Handed back: not in sync with en-US"
$de.Range("J2").Value = "b.md"

$de.Range("A3").Value = "a.md"
$de.Range("H3").Value = "2017-02-21 03:05:14"
$de.Range("J3").Value = "a.md"
$de.Range("L3").Value = "2017-02-21 03:06:52"
$de.Range("M3").Value = "TestHandback_201702211106"

$de.Range("A2").Hyperlinks.Delete()
$de.Hyperlinks.Add($de.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test4/blob/5b8f33dc302b5a2aa99f42855abaa4d3b6b8492e/e2e/a.md", "", "", "b.md") | Out-Null
$de.Hyperlinks.Add($de.Range("J2"), "https://github.com/OpenLocalizationTestOrg/ol-test4-dede/blob/f2aaa14f9f3c3d1cfed34efc5c227b7a115187e2/e2e/a.md", "", "", "b.md") | Out-Null
$de.Hyperlinks.Add($de.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test4/blob/5b8f33dc302b5a2aa99f42855abaa4d3b6b8492e/e2e/b.md", "", "", "a.md") | Out-Null
$de.Hyperlinks.Add($de.Range("J3"), "https://github.com/OpenLocalizationTestOrg/ol-test4-dede/blob/f2aaa14f9f3c3d1cfed34efc5c227b7a115187e2/e2e/b.md", "", "", "a.md") | Out-Null

# ---------------------------------------------------------------------
# Column width tweaks (status / handback-datetime columns grew a bit
# wider to fit the new "not in sync" text).
# ---------------------------------------------------------------------
$ov.Columns.Item(5).ColumnWidth = 32.65
$ov.Columns.Item(6).ColumnWidth = 32.65
$zh.Columns.Item(3).ColumnWidth = 32.65
$de.Columns.Item(3).ColumnWidth = 32.65
